$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the shared header strings: "<Name>_old" -> "<Name>_FV2210" and
# "<Name>_new" -> "<Name>_FV2304" (the "diff" column header is unchanged).
$pairs = @(
  @("Segmentname_old", "Segmentname_FV2210"),
  @("Segmentgruppe_old", "Segmentgruppe_FV2210"),
  @("Segment_old", "Segment_FV2210"),
  @("Datenelement_old", "Datenelement_FV2210"),
  @("Segment ID_old", "Segment ID_FV2210"),
  @("Code_old", "Code_FV2210"),
  @("Qualifier_old", "Qualifier_FV2210"),
  @("Beschreibung_old", "Beschreibung_FV2210"),
  @("Bedingungsausdruck_old", "Bedingungsausdruck_FV2210"),
  @("Bedingung_old", "Bedingung_FV2210"),
  @("Segmentname_new", "Segmentname_FV2304"),
  @("Segmentgruppe_new", "Segmentgruppe_FV2304"),
  @("Segment_new", "Segment_FV2304"),
  @("Datenelement_new", "Datenelement_FV2304"),
  @("Segment ID_new", "Segment ID_FV2304"),
  @("Code_new", "Code_FV2304"),
  @("Qualifier_new", "Qualifier_FV2304"),
  @("Beschreibung_new", "Beschreibung_FV2304"),
  @("Bedingungsausdruck_new", "Bedingungsausdruck_FV2304"),
  @("Bedingung_new", "Bedingung_FV2304")
)

foreach ($p in $pairs) {
  $ws.Cells.Replace($p[0], $p[1], 1, 1, $false, $false, $false, $false)
}

# Turn the data range into a real Excel Table (adds xl/tables/table1.xml,
# the autoFilter, tableParts entry and the worksheet rels).
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U57"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = $null

# Freeze the header row (pane split after row 1).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
